$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 74 (pushes existing rows 74+ down by two)
$ws.Rows.Item(74).Resize(2).Insert()

$ws.Range("A74").Value = "Nestle Library Permanent Reserve"
$ws.Range("E74").Value = "Nestlé Library > Reserve"

$ws.Range("A75").Value = "Nestle Library Reserve"
$ws.Range("E75").Value = "Nestlé Library > Reserve"

$ws.Range("A74").Select()
